$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary block (rows 3-7) ---
$ws.Range("F3").Value = 99648.61
$ws.Range("F4").Value = 1633.801242200003
$ws.Range("F5").Value = 98014.8087578
$ws.Range("F7").Value = 98014.8087578

# --- Row 15 ---
$ws.Range("F15").Value = 381.0163426460692
$ws.Range("G15").Value = 24184.90728980595
$ws.Range("J15").Value = 4837
$ws.Range("K15").Value = 76.2
$ws.Range("L15").Value = 0.133349999999993
$ws.Range("M15").Value = 76.33335
$ws.Range("N15").Value = 4760.8
$ws.Range("O15").Value = 4760.66665
$ws.Range("P15").Value = 4760.66665

# --- Row 16 ---
$ws.Range("F16").Value = 380.8710732876146
$ws.Range("G16").Value = 8212.787507229612
$ws.Range("J16").Value = 9939.940000000001
$ws.Range("K16").Value = 460.97
$ws.Range("L16").Value = 4.973866300000054
$ws.Range("M16").Value = 465.9438663000001
$ws.Range("N16").Value = 9478.970000000001
$ws.Range("O16").Value = 9473.9961337
$ws.Range("P16").Value = 9473.9961337

# --- Row 17 ---
$ws.Range("F17").Value = 380.8619555555555
$ws.Range("G17").Value = 10631.68288888889
$ws.Range("J17").Value = 9568.51
$ws.Range("K17").Value = 342.78
$ws.Range("L17").Value = 0.9460728000000245
$ws.Range("M17").Value = 343.7260728
$ws.Range("N17").Value = 9225.73
$ws.Range("O17").Value = 9224.7839272
$ws.Range("P17").Value = 9224.7839272

# --- Row 18 ---
$ws.Range("F18").Value = 383.67
$ws.Range("K18").Value = 76.73
$ws.Range("L18").Value = 0.7504193999999984
$ws.Range("M18").Value = 77.4804194
$ws.Range("N18").Value = 4259.46
$ws.Range("O18").Value = 4258.7095806
$ws.Range("P18").Value = 4258.7095806

# --- Row 19 ---
$ws.Range("F19").Value = 382.4167058823529
$ws.Range("K19").Value = 65.01000000000001
$ws.Range("L19").Value = 0.4401177000000018
$ws.Range("M19").Value = 65.45011770000001
$ws.Range("N19").Value = 4854.4
$ws.Range("O19").Value = 4853.9598823
$ws.Range("P19").Value = 4853.9598823

# --- Row 20 ---
$ws.Range("F20").Value = 380.9766745453482
$ws.Range("K20").Value = 600.8
$ws.Range("L20").Value = 4.06741599999998
$ws.Range("M20").Value = 604.8674159999999
$ws.Range("N20").Value = 55443.17
$ws.Range("O20").Value = 55439.102584
$ws.Range("P20").Value = 55439.102584
